$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = $val
    $r.Font.Name = "Calibri"
    $r.Font.ThemeColor = 1
    $r.Font.Size = 10
}

# Row 2: MCH188-1
Set-RowCell "A2" "MCH188-1"
Set-RowCell "C2" "CAR CAMPAIGN AGAINST RACIAL EXPLOTATION"
Set-RowCell "D2" ""
Set-RowCell "E2" "Series"
Set-RowCell "F2" "1 Box"
Set-RowCell "G2" "LOCATION: 23O | GRAP COUNT NUMER: NONE"
Set-RowCell "H2" ""

# Row 3: MCH188-2
Set-RowCell "A3" "MCH188-2"
Set-RowCell "C3" "CAR CAMPAIGN AGAINST RACIAL EXPLOTATION"
Set-RowCell "D3" ""
Set-RowCell "E3" "Series"
Set-RowCell "F3" "1 Box"
Set-RowCell "G3" "LOCATION: 23O | GRAP COUNT NUMER: NONE"
Set-RowCell "H3" ""

# Match frozen-pane selection now covering both new rows
$ws.Range("A2:J3").Select()
